$d = $word.ActiveDocument

# --- Locate the final paragraph (contains the last diagram + the _GoBack bookmark) ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)

# Remove the old _GoBack bookmark (it currently sits right before the diagram run).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- Paragraph 1: empty paragraph (NoSpacing, sz 20) ---
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$pr = $p.Range
$pr.Collapse(0)
$pr.Font.Size = 10

# --- Paragraph 2: "TODO:" (NoSpacing, bold, sz 24) ---
$pr.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$pr = $p.Range
$pr.Collapse(0)
$pr.Text = "TODO:"
$pr.Font.Bold = $true
$pr.Font.Size = 12

# --- Paragraph 3: "Write some tests ..." + "moodle" (NoSpacing, sz 20) ---
$pr.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$pr = $p.Range
$pr.Collapse(0)
$pr.Text = [char]0x2013
$pr.Text = "Write some tests " + [char]0x2013 + " Look at the testing video lecture from this week on moodle"
$pr.Font.Bold = $false
$pr.Font.Size = 10

# --- Paragraph 4: "List of things ..." (NoSpacing, sz 20) ---
$pr.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$pr = $p.Range
$pr.Collapse(0)
$pr.Text = "List of things that you have changed from the plan " + [char]0x2013 + " For critical appraisal"
$pr.Font.Size = 10

# Re-create the _GoBack bookmark at the very end of the document (collapsed range).
$endR = $pr.Duplicate
$endR.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endR)
